# Apply edits described by the diff:
# - Set D4:D19 (S_P column, rows 4-19) to 0
# - Change the active selection on the sheet to D20

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set S_P values to 0 for rows 4 through 19 (row 3 is left unchanged)
$ws.Range("D4:D19").Value = 0

# Update the selection to D20 (as seen in the saved worksheet's sheetView)
$ws.Range("D20").Select()
